$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC23_Verify_UserRegistration")

# Insert a new row before row 3 (shifts CLICK/LoginOption etc. down by one)
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with a WAIT step
$ws.Cells.Item(3, 2).Value = "WAIT"

# Update the selection to match the authored state
$ws.Range("C9").Select()
